$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("L3").Value = 0
$ws.Range("N3").ClearContents()
$ws.Range("H3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("M28").Value = 87.85714999999999
$ws.Range("L28").Value = 391
$ws.Range("N28").Value = -1361
$ws.Range("K28").Value = 397.14285
$ws.Range("I28").Value = 397.14285
$ws.Range("H28").Value = 396.73334
$ws.Range("J28").Value = 391
$ws.Range("L87").Value = 25800.56
$ws.Range("N87").Value = -28296.56
$ws.Range("H87").Value = 25800.56
$ws.Range("J87").Value = 25800.56
$ws.Range("L90").Value = 77401.68000000001
$ws.Range("N90").Value = -89881.68000000001
$ws.Range("H90").Value = 25800.56
$ws.Range("J90").Value = 25800.56
$ws.Range("M98").Value = -1148.7273
$ws.Range("L98").Value = 2507.1428
$ws.Range("N98").Value = -5503.1428
$ws.Range("K98").Value = 2646.7273
$ws.Range("I98").Value = 2646.7273
$ws.Range("H98").Value = 2613.0344
$ws.Range("J98").Value = 2507.1428
$ws.Range("L102").Value = 0
$ws.Range("N102").ClearContents()
$ws.Range("H102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("M112").Value = -34892
$ws.Range("L112").Value = 3862.2
$ws.Range("N112").Value = -6078.200000000001
$ws.Range("K112").Value = 36000
$ws.Range("I112").Value = 12000
$ws.Range("H112").Value = 1956.9375
$ws.Range("J112").Value = 1287.4
$ws.Range("M116").Value = 441.4443000000001
$ws.Range("K116").Value = 3000.5557
$ws.Range("I116").Value = 3000.5557
$ws.Range("H116").Value = 3657.353
$ws.Range("M122").Value = -5490.1819
$ws.Range("L122").Value = 7521.428400000001
$ws.Range("N122").Value = -12421.4284
$ws.Range("K122").Value = 7940.1819
$ws.Range("I122").Value = 2646.7273
$ws.Range("H122").Value = 2613.0344
$ws.Range("J122").Value = 2507.1428
$ws.Range("M125").Value = -9960
$ws.Range("L125").Value = 13077
$ws.Range("N125").Value = -17997
$ws.Range("K125").Value = 12420
$ws.Range("I125").Value = 1380
$ws.Range("H125").Value = 1419.8182
$ws.Range("J125").Value = 1453
$ws.Range("L129").Value = 3195.2502
$ws.Range("N129").Value = -13195.2502
$ws.Range("H129").Value = 1191.0741
$ws.Range("J129").Value = 1065.0834

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("M122").Value = -2993.2414
$ws.Range("L122").Value = 10599.2724
$ws.Range("N122").Value = -15499.2724
$ws.Range("K122").Value = 5443.2414
$ws.Range("I122").Value = 1814.4138
$ws.Range("H122").Value = 2287.05
$ws.Range("J122").Value = 3533.0908
$ws.Range("M132").Value = -2820.307699999999
$ws.Range("K132").Value = 5350.307699999999
$ws.Range("I132").Value = 1783.4359
$ws.Range("H132").Value = 2155.52

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("M105").Value = 98.16660000000002
$ws.Range("K105").Value = 1648.8334
$ws.Range("I105").Value = 1648.8334
$ws.Range("H105").Value = 1778.2916

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("M31").Value = -1929.5
$ws.Range("L31").Value = 4324.3687
$ws.Range("N31").Value = -4914.3687
$ws.Range("K31").Value = 2224.5
$ws.Range("I31").Value = 2224.5
$ws.Range("H31").Value = 3006.804
$ws.Range("J31").Value = 4324.3687
$ws.Range("M34").Value = -2022.5
$ws.Range("L34").Value = 4324.3687
$ws.Range("N34").Value = -4728.3687
$ws.Range("K34").Value = 2224.5
$ws.Range("I34").Value = 2224.5
$ws.Range("H34").Value = 3006.804
$ws.Range("J34").Value = 4324.3687
$ws.Range("M105").Value = -690.1819999999998
$ws.Range("L105").Value = 5366.6665
$ws.Range("N105").Value = -8860.666499999999
$ws.Range("K105").Value = 2437.182
$ws.Range("I105").Value = 2437.182
$ws.Range("H105").Value = 3064.9285
$ws.Range("J105").Value = 5366.6665
$ws.Range("M107").Value = 1417.6842
$ws.Range("K107").Value = 502.3158
$ws.Range("I107").Value = 502.3158
$ws.Range("H107").Value = 1379.8636
$ws.Range("M132").Value = -5430.7145
$ws.Range("L132").Value = 10800
$ws.Range("N132").Value = -15860
$ws.Range("K132").Value = 7960.7145
$ws.Range("I132").Value = 2653.5715
$ws.Range("H132").Value = 3070
$ws.Range("J132").Value = 3600

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("M12").Value = 139.499999
$ws.Range("L12").Value = 625.5
$ws.Range("N12").Value = -971.5
$ws.Range("K12").Value = 33.500001
$ws.Range("I12").Value = 11.166667
$ws.Range("H12").Value = 166.21428
$ws.Range("J12").Value = 208.5
$ws.Range("M23").Value = 39.50000800000001
$ws.Range("L23").Value = 369
$ws.Range("N23").Value = -839
$ws.Range("K23").Value = 195.499992
$ws.Range("I23").Value = 65.166664
$ws.Range("H23").Value = 98.21429000000001
$ws.Range("J23").Value = 123
$ws.Range("M81").Value = 286
$ws.Range("L81").Value = 50571.429
$ws.Range("N81").Value = -52817.429
$ws.Range("K81").Value = 837
$ws.Range("I81").Value = 279
$ws.Range("H81").Value = 14784.875
$ws.Range("J81").Value = 16857.143
$ws.Range("M84").Value = 3105
$ws.Range("L84").Value = 151714.287
$ws.Range("N84").Value = -162946.287
$ws.Range("K84").Value = 2511
$ws.Range("I84").Value = 279
$ws.Range("H84").Value = 14784.875
$ws.Range("J84").Value = 16857.143
$ws.Range("M98").Value = 898
$ws.Range("L98").Value = 1748.0001
$ws.Range("N98").Value = -4744.0001
$ws.Range("K98").Value = 600
$ws.Range("I98").Value = 200
$ws.Range("H98").Value = 304.36365
$ws.Range("J98").Value = 582.6667
$ws.Range("M122").Value = -8661.999400000001
$ws.Range("K122").Value = 11111.9994
$ws.Range("I122").Value = 1234.6666
$ws.Range("H122").Value = 2100.4285
$ws.Range("M138").Value = 1967.125
$ws.Range("L138").Value = 5986.0002
$ws.Range("N138").Value = -16266.0002
$ws.Range("K138").Value = 3172.875
$ws.Range("I138").Value = 1057.625
$ws.Range("H138").Value = 1459.5
$ws.Range("J138").Value = 1995.3334

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("M122").Value = -4265.5
$ws.Range("L122").Value = 15187.5
$ws.Range("N122").Value = -20087.5
$ws.Range("K122").Value = 6715.5
$ws.Range("I122").Value = 2238.5
$ws.Range("H122").Value = 3265.4092
$ws.Range("J122").Value = 5062.5
$ws.Range("M126").Value = -11364941
$ws.Range("L126").Value = 11316.7062
$ws.Range("N126").Value = -16256.7062
$ws.Range("K126").Value = 11367411
$ws.Range("I126").Value = 3789137
$ws.Range("H126").Value = 1570130.1
$ws.Range("J126").Value = 3772.2354

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("M7").Value = -3847039.2
$ws.Range("L7").Value = 2039.0625
$ws.Range("N7").Value = -2263.0625
$ws.Range("K7").Value = 3847151.2
$ws.Range("I7").Value = 3847151.2
$ws.Range("H7").Value = 2382346.5
$ws.Range("J7").Value = 2039.0625
$ws.Range("M46").Value = -691.46295
$ws.Range("K46").Value = 879.46295
$ws.Range("I46").Value = 879.46295
$ws.Range("H46").Value = 1233.1904
$ws.Range("L55").Value = 1442.2
$ws.Range("N55").Value = -1788.2
$ws.Range("H55").Value = 1161.6923
$ws.Range("J55").Value = 1442.2
$ws.Range("L74").Value = 40000
$ws.Range("N74").Value = -41996
$ws.Range("H74").Value = 40000
$ws.Range("J74").Value = 40000
$ws.Range("L77").Value = 120000
$ws.Range("N77").Value = -129984
$ws.Range("H77").Value = 40000
$ws.Range("J77").Value = 40000
$ws.Range("M122").Value = -4971.6844
$ws.Range("L122").Value = 12450
$ws.Range("N122").Value = -17350
$ws.Range("K122").Value = 7421.6844
$ws.Range("I122").Value = 2473.8948
$ws.Range("H122").Value = 2765.3914
$ws.Range("J122").Value = 4150
$ws.Range("M126").Value = -11538983.6
$ws.Range("L126").Value = 6117.1875
$ws.Range("N126").Value = -11057.1875
$ws.Range("K126").Value = 11541453.6
$ws.Range("I126").Value = 3847151.2
$ws.Range("H126").Value = 2382346.5
$ws.Range("J126").Value = 2039.0625

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("M122").Value = -1112575
$ws.Range("L122").Value = 11549.625
$ws.Range("N122").Value = -16449.625
$ws.Range("K122").Value = 1115025
$ws.Range("I122").Value = 371675
$ws.Range("H122").Value = 287600.7
$ws.Range("J122").Value = 3849.875
$ws.Range("M132").Value = -1479.1319
$ws.Range("K132").Value = 4009.1319
$ws.Range("I132").Value = 1336.3773
$ws.Range("H132").Value = 3390.672
